{"js": "// Apply the template edits described by the diff:\n// 1. City placeholder text: \"\u0415\u043a\u0430\u0442\u0435\u0440\u0438\u043d\u0431\u0443\u0440\u0433\" -> \"\u041d\u043e\u0432\u043e\u0440\u043e\u0441\u0441\u0438\u0439\u0441\u043a\"\n// 2. Trim the trailing run of spaces left over after \"\u0411\u0418\u041a: 040349602\"\n// 3. Trim the trailing run of spaces left over after \"\u041e\u0413\u0420\u041d: 1172375061891\"\n// 4. Trim the trailing run of spaces left over after \"\u0418\u041d\u041d: 2315996766\"\n// 5. Drop the stray trailing space after the \"\u042d\u043b. \u043f\u043e\u0447\u0442\u0430: @<EMAIL>@\" tag\n//\n// Each edit below locates the (document-unique) original text with\n// `body.search` and rewrites it in place with `insertText(..., replace)`,\n// which keeps the existing run formatting (rFonts/color/etc.) intact.\n\nconst body = context.document.body;\n\nasync function replaceText(needle, replacement) {\n  const results = body.search(needle, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(replacement, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// 1. City name.\nawait replaceText(\"\u0415\u043a\u0430\u0442\u0435\u0440\u0438\u043d\u0431\u0443\u0440\u0433\", \"\u041d\u043e\u0432\u043e\u0440\u043e\u0441\u0441\u0438\u0439\u0441\u043a\");\n\n// 2. \"\u0411\u0418\u041a: 040349602            \" -> \"\u0411\u0418\u041a: 040349602\"\nawait replaceText(\"\u0411\u0418\u041a: 040349602            \", \"\u0411\u0418\u041a: 040349602\");\n\n// 3. \"\u041e\u0413\u0420\u041d: 1172375061891            \" -> \"\u041e\u0413\u0420\u041d: 1172375061891\"\nawait replaceText(\"\u041e\u0413\u0420\u041d: 1172375061891            \", \"\u041e\u0413\u0420\u041d: 1172375061891\");\n\n// 4. \"\u0418\u041d\u041d: 2315996766                                 \" -> \"\u0418\u041d\u041d: 2315996766\"\nawait replaceText(\n  \"\u0418\u041d\u041d: 2315996766                                 \",\n  \"\u0418\u041d\u041d: 2315996766\"\n);\n\n// 5. Drop the trailing space after the email tag in the passport/bank-details paragraph.\nawait replaceText(\"\u042d\u043b. \u043f\u043e\u0447\u0442\u0430: @<EMAIL>@ \", \"\u042d\u043b. \u043f\u043e\u0447\u0442\u0430: @<EMAIL>@\");\n", "ps1": "# Apply the template edits described by the diff:\n# 1. City placeholder text: \"\u0415\u043a\u0430\u0442\u0435\u0440\u0438\u043d\u0431\u0443\u0440\u0433\" -> \"\u041d\u043e\u0432\u043e\u0440\u043e\u0441\u0441\u0438\u0439\u0441\u043a\"\n# 2. Trim the trailing run of spaces left over after \"\u0411\u0418\u041a: 040349602\"\n# 3. Trim the trailing run of spaces left over after \"\u041e\u0413\u0420\u041d: 1172375061891\"\n# 4. Trim the trailing run of spaces left over after \"\u0418\u041d\u041d: 2315996766\"\n# 5. Drop the stray trailing space after the \"\u042d\u043b. \u043f\u043e\u0447\u0442\u0430: @<EMAIL>@\" tag\n#\n# Every edit below is a straight Find/Replace over the whole document\n# (wdReplaceAll) against document-unique text, so each fires exactly once.\n\n$d = $word.ActiveDocument\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n# 1. City name.\n$find = $d.Content.Find\n$find.Execute(\"\u0415\u043a\u0430\u0442\u0435\u0440\u0438\u043d\u0431\u0443\u0440\u0433\", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"\u041d\u043e\u0432\u043e\u0440\u043e\u0441\u0441\u0438\u0439\u0441\u043a\", $wdReplaceAll)\n\n# 2. \"\u0411\u0418\u041a: 040349602            \" -> \"\u0411\u0418\u041a: 040349602\"\n$find = $d.Content.Find\n$find.Execute(\"\u0411\u0418\u041a: 040349602            \", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"\u0411\u0418\u041a: 040349602\", $wdReplaceAll)\n\n# 3. \"\u041e\u0413\u0420\u041d: 1172375061891            \" -> \"\u041e\u0413\u0420\u041d: 1172375061891\"\n$find = $d.Content.Find\n$find.Execute(\"\u041e\u0413\u0420\u041d: 1172375061891            \", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"\u041e\u0413\u0420\u041d: 1172375061891\", $wdReplaceAll)\n\n# 4. \"\u0418\u041d\u041d: 2315996766                                 \" -> \"\u0418\u041d\u041d: 2315996766\"\n$find = $d.Content.Find\n$find.Execute(\"\u0418\u041d\u041d: 2315996766                                 \", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"\u0418\u041d\u041d: 2315996766\", $wdReplaceAll)\n\n# 5. Drop the trailing space after the email tag in the passport/bank-details paragraph.\n$find = $d.Content.Find\n$find.Execute(\"\u042d\u043b. \u043f\u043e\u0447\u0442\u0430: @<EMAIL>@ \", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"\u042d\u043b. \u043f\u043e\u0447\u0442\u0430: @<EMAIL>@\", $wdReplaceAll)\n"}
